$wb = $excel.ActiveWorkbook

# --- Repayment schedule sheet: insert a new column before column N ---
$wsRepay = $wb.Worksheets.Item("Repayment schedule")
$wsRepay.Columns("N:N").Insert()
$wsRepay.Range("N1").ColumnWidth = 10.7109375

# Make "Repayment schedule" the active sheet/tab and set the selection
$wsRepay.Activate()
$wsRepay.Range("S6").Select()
